$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.652425050735474
$ws.Range("B1").Value = 3.919591903686523
$ws.Range("C1").Value = 2.47603702545166
$ws.Range("D1").Value = 0.8201513886451721
$ws.Range("E1").Value = 0.8434818387031555
